# Apply the "new .ttl from Google sheet has been generated" edit to Sheet1.
# Source range grows from A1:T23 to A1:AA22 with many relabeled prefixes,
# metadata values, column headers and concept rows (beer ontology content).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed / new cell values (columns A-T, rows 1-22) ---
$ws.Cells.Item(1, 2).Value = 'http://purl.org/m4m/beer-ontology/'
$ws.Cells.Item(2, 1).Value = 'PREFIX'
$ws.Cells.Item(2, 2).Value = 'beer-onto'
$ws.Cells.Item(2, 3).Value = 'http://purl.org/m4m/beer-ontology/'
$ws.Cells.Item(3, 2).Value = 'pav'
$ws.Cells.Item(3, 3).Value = 'http://purl.org/pav/'
$ws.Cells.Item(4, 2).Value = 'dct'
$ws.Cells.Item(4, 3).Value = 'http://purl.org/dc/terms/'
$ws.Cells.Item(5, 2).Value = 'owl'
$ws.Cells.Item(5, 3).Value = 'http://www.w3.org/2002/07/owl#'
$ws.Cells.Item(6, 2).Value = 'xsd'
$ws.Cells.Item(6, 3).Value = 'http://www.w3.org/2001/XMLSchema#'
$ws.Cells.Item(7, 2).Value = 'skos'
$ws.Cells.Item(7, 3).Value = 'http://www.w3.org/2004/02/skos/core#'
$ws.Cells.Item(8, 2).Value = 'wikidata'
$ws.Cells.Item(8, 3).Value = 'https://www.wikidata.org/wiki/'
$ws.Cells.Item(9, 1).Value = 'dct:title'
$ws.Cells.Item(9, 2).Value = 'Beer ontology'
$ws.Cells.Item(10, 1).Value = 'skos:prefLable'
$ws.Cells.Item(10, 2).Value = 'Beer ontology'
$ws.Cells.Item(11, 2).Value = 'Qualitative and quantitative variables describing beer'
$ws.Cells.Item(12, 2).Value = 'https://www.linkedin.com/in/kristina-tomicic-6bb443108/'
$ws.Cells.Item(17, 1).Value = 'Identifier'
$ws.Cells.Item(17, 2).Value = 'skos:prefLabel@en'
$ws.Cells.Item(17, 3).Value = 'skos:altLabel(separator=",")'
$ws.Cells.Item(17, 4).Value = 'skos:definition@en'
$ws.Cells.Item(17, 5).Value = 'dct:source(separator=",")'
$ws.Cells.Item(17, 6).Value = 'skos:broader(separator=",")'
$ws.Cells.Item(17, 7).Value = 'skos:closeMatch(separator=",")'
$ws.Cells.Item(17, 8).Value = 'skos:exactMatch(separator=",")'
$ws.Cells.Item(17, 9).Value = 'skos:broadMatch(separator=",")'
$ws.Cells.Item(17, 10).Value = 'skos:narrowMatch(separator=",")'
$ws.Cells.Item(17, 11).Value = 'skos:relatedMatch(separator=",")'
$ws.Cells.Item(17, 12).Value = 'owl:deprecated^^xsd:boolean'
$ws.Cells.Item(17, 13).Value = 'dct:isReplacedBy'
$ws.Cells.Item(17, 14).Value = 'skos:editorialNote@en'
$ws.Cells.Item(17, 15).Value = 'dct:creator(separator=",")'
$ws.Cells.Item(17, 16).Value = 'dct:contributor(separator=",")'
$ws.Cells.Item(18, 1).Value = 'beer-onto:Deprecated'
$ws.Cells.Item(18, 2).Value = 'Deprecated'
$ws.Cells.Item(18, 4).Value = 'This is a set of controlled terms which are deprecated.'
$ws.Cells.Item(18, 8).Value = ' '
$ws.Cells.Item(19, 1).Value = 'beer-onto:alc_percentage'
$ws.Cells.Item(19, 2).Value = 'alc_percentage'
$ws.Cells.Item(19, 4).Value = 'Percentage of alcohol in a unit of a beer'
$ws.Cells.Item(20, 1).Value = 'beer-onto:beer_color'
$ws.Cells.Item(20, 2).Value = 'beer_color'
$ws.Cells.Item(20, 4).Value = 'Color shade of a certain beer.'
$ws.Cells.Item(21, 1).Value = 'beer-onto:beer_nutrition'
$ws.Cells.Item(21, 2).Value = 'beer_nutrition'
$ws.Cells.Item(21, 4).Value = 'Calories in a unit of a beer.'
$ws.Cells.Item(22, 1).Value = 'beer-onto:beer_bitterness'
$ws.Cells.Item(22, 2).Value = 'beer_bitterness'
$ws.Cells.Item(22, 4).Value = 'The scale of beer bitterness determining the taste of beer.'

# --- Clear cells whose content was removed, keeping them as empty cells ---
$ws.Cells.Item(1, 3).Value = ''
$ws.Cells.Item(1, 3).NumberFormat = "General"
$ws.Cells.Item(3, 4).Value = ''
$ws.Cells.Item(3, 4).NumberFormat = "General"
$ws.Cells.Item(4, 4).Value = ''
$ws.Cells.Item(4, 4).NumberFormat = "General"
$ws.Cells.Item(5, 4).Value = ''
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(6, 4).Value = ''
$ws.Cells.Item(6, 4).NumberFormat = "General"
$ws.Cells.Item(7, 4).Value = ''
$ws.Cells.Item(7, 4).NumberFormat = "General"
$ws.Cells.Item(8, 4).Value = ''
$ws.Cells.Item(8, 4).NumberFormat = "General"
$ws.Cells.Item(10, 3).Value = ''
$ws.Cells.Item(10, 3).NumberFormat = "General"
$ws.Cells.Item(11, 3).Value = ''
$ws.Cells.Item(11, 3).NumberFormat = "General"
$ws.Cells.Item(12, 3).Value = ''
$ws.Cells.Item(12, 3).NumberFormat = "General"
$ws.Cells.Item(13, 3).Value = ''
$ws.Cells.Item(13, 3).NumberFormat = "General"
$ws.Cells.Item(14, 3).Value = ''
$ws.Cells.Item(14, 3).NumberFormat = "General"
$ws.Cells.Item(15, 3).Value = ''
$ws.Cells.Item(15, 3).NumberFormat = "General"
$ws.Cells.Item(16, 3).Value = ''
$ws.Cells.Item(16, 3).NumberFormat = "General"
$ws.Cells.Item(18, 3).Value = ''
$ws.Cells.Item(18, 3).NumberFormat = "General"
$ws.Cells.Item(18, 5).Value = ''
$ws.Cells.Item(18, 5).NumberFormat = "General"
$ws.Cells.Item(18, 6).Value = ''
$ws.Cells.Item(18, 6).NumberFormat = "General"
$ws.Cells.Item(18, 7).Value = ''
$ws.Cells.Item(18, 7).NumberFormat = "General"
$ws.Cells.Item(18, 9).Value = ''
$ws.Cells.Item(18, 9).NumberFormat = "General"
$ws.Cells.Item(18, 10).Value = ''
$ws.Cells.Item(18, 10).NumberFormat = "General"
$ws.Cells.Item(18, 11).Value = ''
$ws.Cells.Item(18, 11).NumberFormat = "General"
$ws.Cells.Item(18, 12).Value = ''
$ws.Cells.Item(18, 12).NumberFormat = "General"
$ws.Cells.Item(19, 5).Value = ''
$ws.Cells.Item(19, 5).NumberFormat = "General"
$ws.Cells.Item(19, 9).Value = ''
$ws.Cells.Item(19, 9).NumberFormat = "General"
$ws.Cells.Item(20, 3).Value = ''
$ws.Cells.Item(20, 3).NumberFormat = "General"
$ws.Cells.Item(20, 7).Value = ''
$ws.Cells.Item(20, 7).NumberFormat = "General"
$ws.Cells.Item(20, 8).Value = ''
$ws.Cells.Item(20, 8).NumberFormat = "General"
$ws.Cells.Item(20, 11).Value = ''
$ws.Cells.Item(20, 11).NumberFormat = "General"

# --- Remove old row 23 (no longer present in target sheet) ---
$ws.Rows.Item(23).Delete()

# --- Extend used range through column AA for rows 1-22 (new blank cells) ---
for ($r = 1; $r -le 22; $r++) {
    for ($c = 21; $c -le 27; $c++) {
        $ws.Cells.Item($r, $c).NumberFormat = "General"
    }
}

Write-Output ("UsedRange: " + $ws.UsedRange.Address())
